# RPG.docx edits:
#  1. Move the "_GoBack" bookmark from the Dr. Domo "Passive... resistance"
#     paragraph up to the blank paragraph right after "Characters:" (Word
#     bookmarks are unique by name, so adding it at the new spot removes it
#     from the old one automatically).
#  2. Merge the trailing " " + "(NEEDS EDITED ABILITES)" / "(NEEDS ... ABILITY)"
#     runs back together for Sarenssen, Sir Rib Marrow, Isoroku Fleeganschnitzel,
#     and Dr. Domo's header line (visible text itself is unchanged).
#  3. Rename Dr. Domo's ability "Time for surgery" to
#     "Performance enhancing drugs".

$d = $word.ActiveDocument

# --- 1. Relocate the _GoBack bookmark -------------------------------------
$blankPara = $d.Paragraphs(5)
$d.Bookmarks.Add("_GoBack", $blankPara.Range)

# --- 2. Tidy up run splits around the "(NEEDS ...)" notes -----------------
# (Search for just the parenthesized note, so the match starts on the run
#  boundary and only pulls in the short space-only run ahead of it.)
$sarenssen = $d.Paragraphs(14).Range
$sarenssen.Find.Execute("(NEEDS EDITED ABILITES)", $true, $false, $false, $false, $false, `
  $true, 1, $false, "(NEEDS EDITED ABILITES)", 2)

$sirRibMarrow = $d.Paragraphs(32).Range
$sirRibMarrow.Find.Execute("(NEEDS EDITED ABILITES)", $true, $false, $false, $false, $false, `
  $true, 1, $false, "(NEEDS EDITED ABILITES)", 2)

$isoroku = $d.Paragraphs(41).Range
$isoroku.Find.Execute("(NEEDS EDITED ABILITES)", $true, $false, $false, $false, $false, `
  $true, 1, $false, "(NEEDS EDITED ABILITES)", 2)

$drDomoHeader = $d.Paragraphs(49).Range
$drDomoHeader.Find.Execute("(NEEDS ONLY ONE EDITED ABILITY)", $true, $false, $false, $false, $false, `
  $true, 1, $false, "(NEEDS ONLY ONE EDITED ABILITY)", 2)

# --- 3. Rename the ability --------------------------------------------------
$d.Content.Find.Execute("Time for surgery", $true, $false, $false, $false, $false, `
  $true, 1, $false, "Performance enhancing drugs", 2)

Write-Host "Done."
